# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00"); force
# it to be written/stored as literal text so Excel doesn't silently
# reinterpret it as a number, then drop the temporary format again so
# the cell style matches the untouched cells (no explicit 's' attr).
$ws.Range("D2:D51").NumberFormat = "@"

# subscript-3 character used in very small Shiba Inu-style prices (0.0₃0949)
$sub3 = [char]0x2083

$ws.Range("D2").Value = '48.146.85'
$ws.Range("E2").Value = '  +1.82%  '

$ws.Range("D3").Value = '2.508.93'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '321.34'
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").Value = '108.51'
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("E7").Value = '  +0.93%  '

$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +0.84%  '

$ws.Range("D10").Value = '39.83'
$ws.Range("E10").Value = '  +2.05%  '

$ws.Range("D11").Value = '20.17'
$ws.Range("E11").Value = '  +10.01%  '

$ws.Range("D12").Value = '0.0819'
$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("D14").Value = '7.20'
$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("D15").Value = '2.900.14'
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").Value = '2.515.31'
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("E17").Value = '  +0.10%  '

$ws.Range("D18").Value = '47.988.02'
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  +0.56%  '

$ws.Range("D21").Value = "0.0{0}0949" -f $sub3
$ws.Range("E21").Value = '  +1.34%  '

$ws.Range("D22").Value = '2.74'
$ws.Range("E22").Value = '  +0.96%  '

$ws.Range("D23").Value = '72.12'
$ws.Range("E23").Value = '  +2.54%  '

$ws.Range("D24").Value = '277.35'
$ws.Range("E24").Value = '  +13.02%  '

$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  +0.40%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = '25.90'
$ws.Range("E27").Value = '  +0.71%  '

$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  -0.91%  '

$ws.Range("D29").Value = '9.96'
$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("E30").Value = '  +2.65%  '

$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("D32").Value = '49.41'
$ws.Range("E32").Value = '  -0.76%  '

$ws.Range("D33").Value = '19.41'
$ws.Range("E33").Value = '  -3.66%  '

$ws.Range("E34").Value = '  +0.28%  '

$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").Value = '4.60'
$ws.Range("E38").Value = '  -3.31%  '

$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("D40").Value = '122.99'
$ws.Range("E40").Value = '  +4.08%  '

$ws.Range("D41").Value = '0.112'
$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").Value = '2.22'
$ws.Range("E42").Value = '  -0.66%  '

$ws.Range("E43").Value = '  -5.83%  '

$ws.Range("E44").Value = '  +3.43%  '

$ws.Range("D45").Value = '2.000.74'
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").Value = '3.12'
$ws.Range("E46").Value = '  +3.00%  '

$ws.Range("E47").Value = '  +4.09%  '

$ws.Range("E48").Value = '  -0.96%  '

$ws.Range("D49").Value = '9.03'
$ws.Range("E49").Value = '  -1.19%  '

$ws.Range("D50").Value = '5.19'
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("D51").Value = '79.90'
$ws.Range("E51").Value = '  +2.72%  '

# Restore original (default) cell formatting on column D now that the
# text values are committed, so only the displayed text differs.
$ws.Range("D2:D51").ClearFormats()
